$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.999.22'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '3.516.93'
$ws.Range("E3").Value = '  -3.26%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").Value = '3.509.90'
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -4.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.580'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("E14").Value = '  -2.55%  '
$ws.Range("D15").Value = '4.079.98'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '619.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.95%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.062.01'
$ws.Range("E18").Value = '  -2.31%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.513.92'
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.884'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -5.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.60%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.81%  '
$ws.Range("E33").Value = '  -4.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '628.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("E37").Value = '  -12.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.102'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("E39").Value = '  -1.60%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0453'
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("D43").Value = '3.370.74'
$ws.Range("E43").Value = '  -5.11%  '
$ws.Range("E44").Value = '  -5.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("D46").Value = '0.0₃0694'
$ws.Range("E46").Value = '  -5.13%  '
$ws.Range("E47").Value = '  -5.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.129'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.41%  '
$ws.Range("E51").Value = '  +12.79%  '
